$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new headers I1 ("I0") and J1 ("IF")
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy header style (bold, centered, bordered) from H1 onto I1:J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Fill in the I0 / IF numeric data for rows 2-65
$iValues = @(6,6,9,5,8,9,9,8,8,9,9,9,8,8,8,8,8,9,9,7,8,8,8,11,8,8,8,8,6,8,7,7,7,7,6,5,6,7,6,8,9,8,6,8,7,6,1,7,7,7,6,8,7,7,6,4,5,7,8,4,4,5,1,1)
$jValues = @(7,6,9,6,9,9,9,8,9,9,9,9,9,9,9,8,9,9,9,8,8,8,8,11,8,8,8,8,7,8,7,8,7,8,6,7,7,8,6,8,9,8,7,8,8,7,2,7,7,8,8,8,8,7,7,6,7,9,9,5,6,6,3,2)

for ($idx = 0; $idx -lt $iValues.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$idx]
    $ws.Cells.Item($row, 10).Value = $jValues[$idx]
}

